$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "gdsfsd"
$ws.Range("B2").Value = "sdfsdfsdf"

$ws.Range("B2").Select()
